$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Rows.Item(13).Delete()

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = '2024-07-06'
$ws.Range("B2").ClearFormats()
$ws.Range("C2").Value = '南宁·小蜜蜂动漫嘉年华2.0'
$ws.Range("D2").Value = '亭洪路45号 百益上河城'
$ws.Range("E2").Value = '2024.07.06 10:00-07.06 17:00'
$ws.Range("F2").Value = 557
$ws.Range("G2").Value = 50
$ws.Range("H2").Value = 'https://show.bilibili.com/platform/detail.html?id=84925'
$ws.Range("I2").Value = '//i2.hdslb.com/bfs/openplatform/202404/YjFyyYq51713508727131.jpeg'

$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = '2024-07-06'
$ws.Range("B3").ClearFormats()
$ws.Range("C3").Value = '南宁·首届童话梦境Lolita茶会'
$ws.Range("D3").Value = '明秀东路157号 利泰国际大酒店'
$ws.Range("E3").Value = '2024.07.06 13:00-07.06 17:00'
$ws.Range("F3").Value = 183
$ws.Range("G3").Value = 88
$ws.Range("H3").Value = 'https://show.bilibili.com/platform/detail.html?id=85776'
$ws.Range("I3").Value = '//i2.hdslb.com/bfs/openplatform/202405/Xl4NBnky1715847180514.jpeg'

$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = '2024-07-12'
$ws.Range("B4").ClearFormats()
$ws.Range("C4").Value = '南宁·漫控嘉年华09暨南宁高校动漫联盟十六周年联合漫展'
$ws.Range("D4").Value = '民族大道106号 南宁国际会展中心'
$ws.Range("E4").Value = '2024.07.12 09:30-07.14 17:00'
$ws.Range("F4").Value = 314
$ws.Range("G4").Value = 50
$ws.Range("H4").Value = 'https://show.bilibili.com/platform/detail.html?id=87182'
$ws.Range("I4").Value = '//i1.hdslb.com/bfs/openplatform/202406/x4UZPn301718159475475.jpeg'

$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = '2024-07-13'
$ws.Range("B5").ClearFormats()
$ws.Range("C5").Value = '南宁·0713国乙ONLY'
$ws.Range("D5").Value = '亭洪路45号 水明漾宴会中心'
$ws.Range("E5").Value = '2024.07.13 09:30-07.13 21:00'
$ws.Range("F5").Value = 400
$ws.Range("G5").Value = 68
$ws.Range("H5").Value = 'https://show.bilibili.com/platform/detail.html?id=86378'
$ws.Range("I5").Value = '//i1.hdslb.com/bfs/openplatform/202405/ZDBCv2of1716659486569.jpeg'

$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = '2024-07-14'
$ws.Range("B6").ClearFormats()
$ws.Range("C6").Value = '广西·首届明日方舟only展 - 花庭圣梦'
$ws.Range("D6").Value = '明秀东路157号 利泰国际大酒店'
$ws.Range("E6").Value = '2024.07.14 09:00-07.14 18:00'
$ws.Range("F6").Value = 255
$ws.Range("G6").Value = 69
$ws.Range("H6").Value = 'https://show.bilibili.com/platform/detail.html?id=85852'
$ws.Range("I6").Value = '//i2.hdslb.com/bfs/openplatform/202405/xsMTmueN1715920435584.jpeg'

$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = '2024-07-20'
$ws.Range("B7").ClearFormats()
$ws.Range("C7").Value = '南宁·AB动漫游戏嘉年华'
$ws.Range("D7").Value = '三塘南路与长虹东路交叉路口往北约50米 广西农业会展中心'
$ws.Range("E7").Value = '2024.07.20 09:30-07.21 17:00'
$ws.Range("F7").Value = 2354
$ws.Range("G7").Value = 60
$ws.Range("H7").Value = 'https://show.bilibili.com/platform/detail.html?id=84862'
$ws.Range("I7").Value = '//i1.hdslb.com/bfs/openplatform/202404/eglavDeZ1714036487217.jpeg'

$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = '2024-07-20'
$ws.Range("B8").ClearFormats()
$ws.Range("C8").Value = '横州·第二届海棠动漫游戏嘉年华'
$ws.Range("D8").Value = '茉莉花大道 横州国际大酒店'
$ws.Range("E8").Value = '2024.07.20 09:30-07.20 17:00'
$ws.Range("F8").Value = 392
$ws.Range("G8").Value = 30
$ws.Range("H8").Value = 'https://show.bilibili.com/platform/detail.html?id=84799'
$ws.Range("I8").Value = '//i2.hdslb.com/bfs/openplatform/202404/r50S2ttT1713869164413.jpeg'

$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = '2024-07-27'
$ws.Range("B9").ClearFormats()
$ws.Range("C9").Value = '南宁·第十九届（2024）良牙动漫夏季盛典（良牙夏典）'
$ws.Range("D9").Value = '民族大道106号 南宁国际会展中心'
$ws.Range("E9").Value = '2024.07.27 09:30-07.28 17:30'
$ws.Range("F9").Value = 5977
$ws.Range("G9").Value = 55
$ws.Range("H9").Value = 'https://show.bilibili.com/platform/detail.html?id=85264'
$ws.Range("I9").Value = '//i1.hdslb.com/bfs/openplatform/202406/JxFed5iv1718622152091.jpeg'

$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = '2024-08-03'
$ws.Range("B10").ClearFormats()
$ws.Range("C10").Value = '南宁·火影忍者only'
$ws.Range("D10").Value = '厢竹大道65号 桔子酒店'
$ws.Range("E10").Value = '2024.08.03 10:00-08.03 17:00'
$ws.Range("F10").Value = 151
$ws.Range("G10").Value = 68
$ws.Range("H10").Value = 'https://show.bilibili.com/platform/detail.html?id=86994'
$ws.Range("I10").Value = '//i0.hdslb.com/bfs/openplatform/202406/h1tXE9t11717523356034.jpeg'

$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = '2024-08-03'
$ws.Range("B11").ClearFormats()
$ws.Range("C11").Value = '南宁·蔚蓝档案only'
$ws.Range("D11").Value = '亭洪路45号 百益上河城'
$ws.Range("E11").Value = '2024.08.03 09:00-08.03 17:00'
$ws.Range("F11").Value = 383
$ws.Range("G11").Value = 68
$ws.Range("H11").Value = 'https://show.bilibili.com/platform/detail.html?id=85370'
$ws.Range("I11").Value = '//i1.hdslb.com/bfs/openplatform/202405/sBxi2Mx61715247424836.jpeg'

$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = '2024-11-02'
$ws.Range("B12").ClearFormats()
$ws.Range("C12").Value = '南宁·万圣漫控嘉年华10'
$ws.Range("D12").Value = '亭洪路45号 百益上河城'
$ws.Range("E12").Value = '2024.11.02 11:00-11.03 22:00'
$ws.Range("F12").Value = 12
$ws.Range("G12").Value = 50
$ws.Range("H12").Value = 'https://show.bilibili.com/platform/detail.html?id=87820'
$ws.Range("I12").Value = '//i1.hdslb.com/bfs/openplatform/202406/abJD2cvV1718955681653.jpeg'

$ws = $wb.Worksheets.Item(4)
$ws.Rows.Item(17).Delete()

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = '2024-06-22'
$ws.Range("B2").ClearFormats()
$ws.Range("C2").Value = '南宁·浪漫古典·百年经典世界名曲音乐会'
$ws.Range("D2").Value = '广西壮族自治区南宁市良庆区龙堤路25号  广西文化艺术中心-音乐厅'
$ws.Range("E2").Value = '2024.06.22 20:00-06.22 21:30'
$ws.Range("F2").Value = 52
$ws.Range("G2").Value = 135
$ws.Range("H2").Value = 'https://show.bilibili.com/platform/detail.html?id=83959'
$ws.Range("I2").Value = '//i1.hdslb.com/bfs/openplatform/202404/H0f8U7no1712041461015.jpeg'

$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = '2024-07-06'
$ws.Range("B3").ClearFormats()
$ws.Range("C3").Value = '南宁·小蜜蜂动漫嘉年华2.0'
$ws.Range("D3").Value = '亭洪路45号 百益上河城'
$ws.Range("E3").Value = '2024.07.06 10:00-07.06 17:00'
$ws.Range("F3").Value = 557
$ws.Range("G3").Value = 50
$ws.Range("H3").Value = 'https://show.bilibili.com/platform/detail.html?id=84925'
$ws.Range("I3").Value = '//i2.hdslb.com/bfs/openplatform/202404/YjFyyYq51713508727131.jpeg'

$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = '2024-07-06'
$ws.Range("B4").ClearFormats()
$ws.Range("C4").Value = '南宁·首届童话梦境Lolita茶会'
$ws.Range("D4").Value = '明秀东路157号 利泰国际大酒店'
$ws.Range("E4").Value = '2024.07.06 13:00-07.06 17:00'
$ws.Range("F4").Value = 183
$ws.Range("G4").Value = 88
$ws.Range("H4").Value = 'https://show.bilibili.com/platform/detail.html?id=85776'
$ws.Range("I4").Value = '//i2.hdslb.com/bfs/openplatform/202405/Xl4NBnky1715847180514.jpeg'

$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = '2024-07-12'
$ws.Range("B5").ClearFormats()
$ws.Range("C5").Value = '南宁·漫控嘉年华09暨南宁高校动漫联盟十六周年联合漫展'
$ws.Range("D5").Value = '民族大道106号 南宁国际会展中心'
$ws.Range("E5").Value = '2024.07.12 09:30-07.14 17:00'
$ws.Range("F5").Value = 314
$ws.Range("G5").Value = 50
$ws.Range("H5").Value = 'https://show.bilibili.com/platform/detail.html?id=87182'
$ws.Range("I5").Value = '//i1.hdslb.com/bfs/openplatform/202406/x4UZPn301718159475475.jpeg'

$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = '2024-07-13'
$ws.Range("B6").ClearFormats()
$ws.Range("C6").Value = '南宁·0713国乙ONLY'
$ws.Range("D6").Value = '亭洪路45号 水明漾宴会中心'
$ws.Range("E6").Value = '2024.07.13 09:30-07.13 21:00'
$ws.Range("F6").Value = 400
$ws.Range("G6").Value = 68
$ws.Range("H6").Value = 'https://show.bilibili.com/platform/detail.html?id=86378'
$ws.Range("I6").Value = '//i1.hdslb.com/bfs/openplatform/202405/ZDBCv2of1716659486569.jpeg'

$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = '2024-07-14'
$ws.Range("B7").ClearFormats()
$ws.Range("C7").Value = '广西·首届明日方舟only展 - 花庭圣梦'
$ws.Range("D7").Value = '明秀东路157号 利泰国际大酒店'
$ws.Range("E7").Value = '2024.07.14 09:00-07.14 18:00'
$ws.Range("F7").Value = 255
$ws.Range("G7").Value = 69
$ws.Range("H7").Value = 'https://show.bilibili.com/platform/detail.html?id=85852'
$ws.Range("I7").Value = '//i2.hdslb.com/bfs/openplatform/202405/xsMTmueN1715920435584.jpeg'

$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = '2024-07-18'
$ws.Range("B8").ClearFormats()
$ws.Range("C8").Value = '南宁·限时6折|俄罗斯圣彼得堡古典芭蕾舞剧院《胡桃夹子》'
$ws.Range("D8").Value = '龙堤路25号 广西文化艺术中心'
$ws.Range("E8").Value = '2024.07.18 20:00-07.18 21:30'
$ws.Range("F8").Value = 11
$ws.Range("G8").Value = 108
$ws.Range("H8").Value = 'https://show.bilibili.com/platform/detail.html?id=85816'
$ws.Range("I8").Value = '//i0.hdslb.com/bfs/openplatform/202405/SN0ZyGVj1715675672714.jpeg'

$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = '2024-07-19'
$ws.Range("B9").ClearFormats()
$ws.Range("C9").Value = '南宁·限时6折|俄罗斯圣彼得堡古典芭蕾舞剧院《天鹅湖》 '
$ws.Range("D9").Value = '龙堤路25号 广西文化艺术中心'
$ws.Range("E9").Value = '2024.07.19 20:00-07.19 22:00'
$ws.Range("F9").Value = 15
$ws.Range("G9").Value = 108
$ws.Range("H9").Value = 'https://show.bilibili.com/platform/detail.html?id=85831'
$ws.Range("I9").Value = '//i1.hdslb.com/bfs/openplatform/202405/ZyyeeOUo1715677877362.jpeg'

$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = '2024-07-20'
$ws.Range("B10").ClearFormats()
$ws.Range("C10").Value = '南宁·AB动漫游戏嘉年华'
$ws.Range("D10").Value = '三塘南路与长虹东路交叉路口往北约50米 广西农业会展中心'
$ws.Range("E10").Value = '2024.07.20 09:30-07.21 17:00'
$ws.Range("F10").Value = 2354
$ws.Range("G10").Value = 60
$ws.Range("H10").Value = 'https://show.bilibili.com/platform/detail.html?id=84862'
$ws.Range("I10").Value = '//i1.hdslb.com/bfs/openplatform/202404/eglavDeZ1714036487217.jpeg'

$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = '2024-07-20'
$ws.Range("B11").ClearFormats()
$ws.Range("C11").Value = '横州·第二届海棠动漫游戏嘉年华'
$ws.Range("D11").Value = '茉莉花大道 横州国际大酒店'
$ws.Range("E11").Value = '2024.07.20 09:30-07.20 17:00'
$ws.Range("F11").Value = 392
$ws.Range("G11").Value = 30
$ws.Range("H11").Value = 'https://show.bilibili.com/platform/detail.html?id=84799'
$ws.Range("I11").Value = '//i2.hdslb.com/bfs/openplatform/202404/r50S2ttT1713869164413.jpeg'

$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = '2024-07-27'
$ws.Range("B12").ClearFormats()
$ws.Range("C12").Value = '南宁·第十九届（2024）良牙动漫夏季盛典（良牙夏典）'
$ws.Range("D12").Value = '民族大道106号 南宁国际会展中心'
$ws.Range("E12").Value = '2024.07.27 09:30-07.28 17:30'
$ws.Range("F12").Value = 5977
$ws.Range("G12").Value = 55
$ws.Range("H12").Value = 'https://show.bilibili.com/platform/detail.html?id=85264'
$ws.Range("I12").Value = '//i1.hdslb.com/bfs/openplatform/202406/JxFed5iv1718622152091.jpeg'

$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = '2024-08-03'
$ws.Range("B13").ClearFormats()
$ws.Range("C13").Value = '南宁·火影忍者only'
$ws.Range("D13").Value = '厢竹大道65号 桔子酒店'
$ws.Range("E13").Value = '2024.08.03 10:00-08.03 17:00'
$ws.Range("F13").Value = 151
$ws.Range("G13").Value = 68
$ws.Range("H13").Value = 'https://show.bilibili.com/platform/detail.html?id=86994'
$ws.Range("I13").Value = '//i0.hdslb.com/bfs/openplatform/202406/h1tXE9t11717523356034.jpeg'

$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = '2024-08-03'
$ws.Range("B14").ClearFormats()
$ws.Range("C14").Value = '南宁·蔚蓝档案only'
$ws.Range("D14").Value = '亭洪路45号 百益上河城'
$ws.Range("E14").Value = '2024.08.03 09:00-08.03 17:00'
$ws.Range("F14").Value = 383
$ws.Range("G14").Value = 68
$ws.Range("H14").Value = 'https://show.bilibili.com/platform/detail.html?id=85370'
$ws.Range("I14").Value = '//i1.hdslb.com/bfs/openplatform/202405/sBxi2Mx61715247424836.jpeg'

$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = '2024-08-10'
$ws.Range("B15").ClearFormats()
$ws.Range("C15").Value = '南宁·限时7折|浪漫七夕《一生所爱》《爱乐之城》《假如爱有天意》经典浪漫电影主题音乐会'
$ws.Range("D15").Value = '龙堤路25号 广西文化艺术中心'
$ws.Range("E15").Value = '2024.08.10 20:00-08.10 21:30'
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 99
$ws.Range("H15").Value = 'https://show.bilibili.com/platform/detail.html?id=87729'
$ws.Range("I15").Value = '//i1.hdslb.com/bfs/openplatform/202406/qKUDMYOh1718177639735.png'

$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = '2024-11-02'
$ws.Range("B16").ClearFormats()
$ws.Range("C16").Value = '南宁·万圣漫控嘉年华10'
$ws.Range("D16").Value = '亭洪路45号 百益上河城'
$ws.Range("E16").Value = '2024.11.02 11:00-11.03 22:00'
$ws.Range("F16").Value = 12
$ws.Range("G16").Value = 50
$ws.Range("H16").Value = 'https://show.bilibili.com/platform/detail.html?id=87820'
$ws.Range("I16").Value = '//i1.hdslb.com/bfs/openplatform/202406/abJD2cvV1718955681653.jpeg'
